$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "[-, -, 'MCT-3A-Processos de Usinagem 1', -]"
$ws.Range("E4").Value = "[-, -, 'MCT-3A-Processos de Usinagem 1', -]"
$ws.Range("E6").Value = "[-, -, 'MCT-3A-Processos de Usinagem 1', -]"
$ws.Range("E7").Value = "[-, -, 'MCT-3A-Processos de Usinagem 1', -]"
$ws.Range("E8").Value = "-"
